$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Chief Financial Officer"
$ws.Range("C3").Value = "Carolina Walther"
$ws.Range("D3").Value = "carolina@sesamy.com"
